$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "is_permanent" column (K) with its header.
$ws.Range("K1").Value = "是否常驻。常驻则填 yes"

# Mark the first item (木头/log) as permanent.
$ws.Range("K2").Value = "yes"

# Match the header cell's look (fill/border/alignment/wrap) to the rest of row 1.
$ws.Range("K1").Interior.ThemeColor = $ws.Range("J1").Interior.ThemeColor
$ws.Range("K1").Interior.TintAndShade = $ws.Range("J1").Interior.TintAndShade
$ws.Range("K1").Interior.Pattern = $ws.Range("J1").Interior.Pattern
$ws.Range("K1").Borders.LineStyle = $ws.Range("J1").Borders.LineStyle
$ws.Range("K1").HorizontalAlignment = $ws.Range("J1").HorizontalAlignment
$ws.Range("K1").VerticalAlignment = $ws.Range("J1").VerticalAlignment
$ws.Range("K1").WrapText = $ws.Range("J1").WrapText

# Match the data cell's look (border/alignment) to the rest of row 2.
$ws.Range("K2").Borders.LineStyle = $ws.Range("I2").Borders.LineStyle
$ws.Range("K2").HorizontalAlignment = $ws.Range("I2").HorizontalAlignment
$ws.Range("K2").VerticalAlignment = $ws.Range("I2").VerticalAlignment

# Move/keep the active selection on the newly added column, like the source edit did.
$ws.Range("K3").Select() | Out-Null
